$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting existing rows 7-23 down to 8-24.
$ws.Rows.Item(7).Insert()

# Populate the new weekly record in row 7.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44575
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112045
$ws.Range("G7").Value = "Zapallo"
$ws.Range("H7").Value = "Camote"
$ws.Range("I7").Value = "1a nueva(o)"
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 550
$ws.Range("M7").Value = 525
$ws.Range("N7").Value = '$/kilo (volumen en unidades)'
$ws.Range("O7").Value = "Región de O'Higgins"
$ws.Range("P7").Value = 525
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
